# Update "想去人数" (number of people wishing to attend) values on the
# "展览" and "全部类型" worksheets for rows 2, 3, 5, 6 in column F.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 5692
    $ws.Range("F3").Value = 8
    $ws.Range("F5").Value = 961
    $ws.Range("F6").Value = 44
}
